# Weekly Fruta/Hortaliza update: a new observation is inserted as row 4
# (Brooks, "Segunda" quality, serial date 44571 = 2022-01-10), pushing the
# existing rows 4-10 down to rows 5-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4; this shifts old rows 4..10
# down to 5..11 (values & formatting move with them).
$ws.Rows(4).Insert()

# Populate the newly inserted row 4 with the new weekly record.
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44571
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100103
$ws.Range("H4").Value = "Frutos de hueso (carozo)"
$ws.Range("I4").Value = 100103001
$ws.Range("J4").Value = "Cereza"
$ws.Range("K4").Value = "Brooks"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 400
$ws.Range("N4").Value = 8500
$ws.Range("O4").Value = 9000
$ws.Range("P4").Value = 8750
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 875
$ws.Range("T4").Value = 10
